# 16 Jan 2024 Update - Added new variables to codebook
#
# Adds two new rows to the "Codebook" sheet describing the "Age" and
# "Political party" variables, and leaves the "Codebook" sheet as the
# active/selected sheet (instead of "Data").

$wb = $excel.ActiveWorkbook
$wsCodebook = $wb.Worksheets.Item("Codebook")

# New codebook rows (row 5: Age, row 6: Political party)
$wsCodebook.Range("A5").Value = "Age"
$wsCodebook.Range("B5").Value = "Age in years"
$wsCodebook.Range("C5").Value = "numeric value >0 or NA"

$wsCodebook.Range("A6").Value = "Political party"
$wsCodebook.Range("B6").Value = "Political party affiliation (categorical)"
$wsCodebook.Range("C6").Value = "Republican/Democrat/Libertarian/Green/NA"

# Make "Codebook" the active sheet/tab and move its selection to C7,
# matching the author's view state when they saved the file.
$wsCodebook.Activate()
$wsCodebook.Range("C7").Select() | Out-Null
